$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.416.93"
$ws.Range("E2").Value = "  +1.59%  "
$ws.Range("D3").Value = "3.847.56"
$ws.Range("E3").Value = "  +1.72%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "449.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +7.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.12"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +14.98%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.627"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.89%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.746"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.159"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000327"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.75%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.95"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +10.68%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.45"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.56%  "
$ws.Range("D14").Value = "4.467.93"
$ws.Range("E14").Value = "  +1.86%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.10"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.29%  "
$ws.Range("D16").Value = "3.841.88"
$ws.Range("E16").Value = "  +1.55%  "
$ws.Range("E17").Value = "  -0.20%  "
$ws.Range("E18").Value = "  +3.78%  "
$ws.Range("E19").Value = "  +8.70%  "
$ws.Range("D20").Value = "67.468.57"
$ws.Range("E20").Value = "  +1.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "427.82"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.68%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.79"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.64%  "
$ws.Range("E23").Value = "  +9.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "86.67"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "37.60"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.84%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.45"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +9.05%  "
$ws.Range("E27").Value = "  -2.18%  "
$ws.Range("B28").Value = "Filecoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.80"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.61%  "
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.53"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +18.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "747.47"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.79"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +13.06%  "
$ws.Range("E32").Value = "  +12.53%  "
$ws.Range("E33").Value = "  +1.54%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "43.37"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +14.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.156"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.90%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "57.65"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.82%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.52"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +17.00%  "
$ws.Range("E39").Value = "  +6.25%  "
$ws.Range("B40").Value = "PEPE"
$ws.Range("C40").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D40").Value = "0.0₃0694"
$ws.Range("E40").Value = "  -8.71%  "
$ws.Range("B41").Value = "ThetaToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.92"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.57%  "
$ws.Range("E42").Value = "  +17.94%  "
$ws.Range("E43").Value = "  +5.65%  "
$ws.Range("E44").Value = "  -0.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.48"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.81%  "
$ws.Range("B46").Value = "ARBITRUM"
$ws.Range("C46").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.16"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.95%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.23"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.47"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +12.34%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "146.76"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.66%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.66"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.83%  "
$ws.Range("E51").Value = "  +5.87%  "
